# Append the new Nalco run-log row (row 30) that records the
# 2025-08-19 03:52:09 UTC "SKIPPED" run, matching the formatting of the
# preceding log rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (29) down onto the new
# row (30) so the appended row keeps the same style (s="3") as every other
# log row.
$ws.Range("A29:H29").Copy()
$ws.Range("A30:H30").PasteSpecial(-4122)

# Fill in the new run's data.
$ws.Range("A30").Value = "2025-08-19 03:52:09 UTC"
$ws.Range("B30").Value = "2025-08-19 09:22:09 IST"
$ws.Range("C30").Value = "SKIPPED"
$ws.Range("D30").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E30").Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Range("G30").Value = 0

Write-Output "Appended row 30 to $($ws.Name)"
